$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r Workspace_BusinessLine")

# New rows to append, following the same pattern as the existing row 3
# (Action, Id, Name, <blank Description>, Variable, Business_Line)
$codes = @("000002", "000003", "000004", "000005", "000006")

$r = 4
foreach ($code in $codes) {
    $ws.Cells.Item($r, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($r, 5).Value = "WS_ISPRO_CORPORATE"
    $r = $r + 1
}

$r = 4
foreach ($code in $codes) {
    $blCode = "BL_ISPRO_$code"
    $ws.Cells.Item($r, 6).Value = $blCode
    $r = $r + 1
}

$r = 4
foreach ($code in $codes) {
    $wsCode = "WS_ISPRO_BL_ISPRO_$code"
    $ws.Cells.Item($r, 2).Value = $wsCode
    $ws.Cells.Item($r, 3).Value = $wsCode
    $r = $r + 1
}

# Re-fit the columns to the new content (widths recalculated by Excel's
# "best fit" column-width logic after the new rows were added)
$ws.Columns.Item(1).ColumnWidth = 17.8329
$ws.Columns.Item(2).ColumnWidth = 26.166
$ws.Columns.Item(3).ColumnWidth = 26.166
$ws.Columns.Item(4).ColumnWidth = 10.166
$ws.Columns.Item(6).ColumnWidth = 15.5858

# Make this sheet the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("L11").Select() | Out-Null
